$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.505.28"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "3.450.14"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'579.68"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "'148.23"
$ws.Range("E6").Value = "  +8.59%  "
$ws.Range("D7").Value = "3.450.78"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("E10").Value = "  +3.76%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "'0.391"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "4.040.87"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").Value = "'28.01"
$ws.Range("E14").Value = "  +6.00%  "
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "3.452.19"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "61.624.12"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = "  +8.32%  "
$ws.Range("D20").Value = "'14.34"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").Value = "'9.45"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "'386.33"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").Value = "'0.570"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("D24").Value = "3.597.77"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("D25").Value = "'72.70"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  +7.45%  "
$ws.Range("D30").Value = "'7.83"
$ws.Range("E30").Value = "  +3.79%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'1.53"
$ws.Range("E32").Value = "  -13.58%  "
$ws.Range("D33").Value = "'8.26"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'23.98"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("D38").Value = "'5.24"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").Value = "'166.47"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "'0.0790"
$ws.Range("E41").Value = "  +4.83%  "
$ws.Range("D42").Value = "'26.17"
$ws.Range("E42").Value = "  +9.15%  "
$ws.Range("D43").Value = "'0.797"
$ws.Range("E43").Value = "  +3.03%  "
$ws.Range("D44").Value = "'4.52"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'42.31"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").Value = "2.609.46"
$ws.Range("E48").Value = "  +10.07%  "
$ws.Range("E49").Value = "  -4.05%  "
$ws.Range("D50").Value = "'6.98"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "'23.28"
$ws.Range("E51").Value = "  -0.79%  "
